$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(-37.713, -16.1397),
    @(-26.0908, -16.8939),
    @(-34.6688, -28.9877),
    @(-16.0498, -15.5363),
    @(-9.01318, -17.9497),
    @(1.5022, -15.6872),
    @(6.56223, -16.1397),
    @(-24.3514, -28.3575),
    @(-15.1801, -28.8101),
    @(-7.9063, -28.3575),
    @(6.95754, -27.4525),
    @(8.97106, -32.3006)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$null = $ws.Range("A1:B12").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
